# Regenerate orders with updated distance/size codes.
# The stimulus naming scheme encodes Distance (D##) and Size (S##) tokens
# inside Condition / Filename_Left / Filename_Right / Distance / Size
# columns. This run's distances/sizes shifted:
#   D51 -> D55, D64 -> D69, D80 -> D86, S30 -> S31
# Apply as a whole-workbook text substitution (mirrors how the source
# generator re-stamped every affected string), so every cell/shared string
# built from these tokens (e.g. "Face02_D51_S30", "Face02_D51_S30_l.png",
# "D51", "S30") is updated consistently.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = $ws.Cells

# Excel constants (spelled out explicitly rather than relying on
# auto-$null-ing unset variables):
#   xlPart    = 2  (LookAt: match part of the cell text)
#   xlByRows  = 1  (SearchOrder: irrelevant for a full replace-all, but
#                   must be a valid value)
$xlPart = 2
$xlByRows = 1

# Order doesn't matter here: none of the replacement targets (D55, D69,
# D86, S31) collide with any of the other replacement sources (D51, D64,
# D80, S30), so a straightforward sequential partial-text replace across
# every cell is safe and idempotent.
$cells.Replace("D51", "D55", $xlPart, $xlByRows, $false, $false, $false, $false)
$cells.Replace("D64", "D69", $xlPart, $xlByRows, $false, $false, $false, $false)
$cells.Replace("D80", "D86", $xlPart, $xlByRows, $false, $false, $false, $false)
$cells.Replace("S30", "S31", $xlPart, $xlByRows, $false, $false, $false, $false)

Write-Output "replace-done"
